# "Se agrega reporte de zapatillas"
#
# The task list on Hoja1 drops the row "vista detalle de cuota (venta
# asociada)" (old row 13) - everything below it shifts up one row - the
# "reporte de zapatilla..." task's status cell is reset from "en proceso"
# to a 100% (done) checkbox, and a brand new task "revisar reportes filtro
# fecha" is appended as the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the "vista detalle de cuota (venta asociada)" row entirely; every
# row below (14..40) shifts up to (13..39) automatically.
$ws.Rows.Item(13).Delete()

# The "reporte de zapatilla por talle en stock - historial de compras" task
# (now row 23, previously row 24) had its estado cell switched from the
# text "en proceso" to the same 0%-style numeric "done" flag used by the
# other rows.
$ws.Range("C23").Value = 1
$ws.Range("C23").NumberFormat = "0%"

# New task appended at the bottom of the list.
$ws.Range("A40").Value = "revisar reportes filtro fecha"

# Update the view state to match: scrolled down a bit further, and the
# active selection resting just past the new last row.
$ws.Select()
try {
    $excel.ActiveWindow.ScrollRow = 19
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("A41").Select()
